$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text correction: "Forecasting" -> "Forecasting or Inference" (row 19) ---
$ws.Range("C19").Value = "Forecasting or Inference"

# --- New column L: "Ecology" formatted values (roughly F+I totals) ---

# Row 1: date value, matching I1's date style (copy format so no new numFmt/style is minted)
$ws.Range("L1").Value = 41610
$ws.Range("I1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = 41610

# Row 2: header label "Pages" (same shared text as F2/I2)
$ws.Range("L2").Value = "Pages"

# Row 3: total
$ws.Range("L3").Value = 30

# Section value rows
$ws.Range("L4").Value = 1.5
$ws.Range("L5").Value = 3.5
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 0.7
$ws.Range("L9").Value = 0.9
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L13").Value = 5.15
$ws.Range("L14").Value = 0.95
$ws.Range("L15").Value = 0.75
$ws.Range("L16").Value = 0.6
$ws.Range("L17").Value = 1.1
$ws.Range("L18").Value = 0.6
$ws.Range("L19").Value = 0.7
$ws.Range("L20").Value = 2.1
$ws.Range("L21").Value = 0.7
$ws.Range("L22").Value = 1.3
$ws.Range("L23").Value = 1.7
$ws.Range("L24").Value = 0.6
$ws.Range("L25").Value = 1.05
$ws.Range("L26").Value = 4
$ws.Range("L27").Value = 0.15
$ws.Range("L28").Value = 4.75
$ws.Range("L29").Value = 1
$ws.Range("L30").Value = 3

# Row 32: sum formula of the major section subtotals
$ws.Range("L32").Formula = "=L4+L5+L6+L13+L20+L23+L26+L27+L28+L29+L30"

# Column L width to roughly match column I's best-fit width
$ws.Columns.Item(12).ColumnWidth = 8.83

# Selection ends on the new date cell, like the source workbook
$ws.Range("L1").Select() | Out-Null
